$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.122774
$ws.Range("H2").Value = 33.368322
$ws.Range("I2").Value = 0.2449652610853511
$ws.Range("J2").Value = 0.2449652610853511
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.127188333333333
$ws.Range("N2").Value = 12.381565
$ws.Range("O2").Value = 0.0561359176022362
$ws.Range("P2").Value = 0.05613591760223619
$ws.Range("Q2").Value = 45.90578308710332
$ws.Range("R2").Value = 413.15204778393
$ws.Range("S2").Value = 0.01375134971169755
$ws.Range("T2").Value = 0.01375134971169755
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.122774
$ws.Range("H3").Value = 33.368322
$ws.Range("I3").Value = 0.2449652610853511
$ws.Range("J3").Value = 0.2449652610853511
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.24901333333333
$ws.Range("N3").Value = 141.74704
$ws.Range("O3").Value = 0.6426570597336346
$ws.Range("P3").Value = 0.6426570597336345
$ws.Range("Q3").Value = 525.5400970296533
$ws.Range("R3").Value = 4729.86087326688
$ws.Range("S3").Value = 0.1574286544259939
$ws.Range("T3").Value = 0.1574286544259939
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.122774
$ws.Range("H4").Value = 33.368322
$ws.Range("I4").Value = 0.2449652610853511
$ws.Range("J4").Value = 0.2449652610853511
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.218847999999999
$ws.Range("N4").Value = 24.656544
$ws.Range("O4").Value = 0.1117885923419141
$ws.Range("P4").Value = 0.1117885923419141
$ws.Range("Q4").Value = 91.41638884435199
$ws.Range("R4").Value = 822.747499599168
$ws.Range("S4").Value = 0.02738432170940086
$ws.Range("T4").Value = 0.02738432170940086
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.122774
$ws.Range("H5").Value = 33.368322
$ws.Range("I5").Value = 0.2449652610853511
$ws.Range("J5").Value = 0.2449652610853511
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.92629833333334
$ws.Range("N5").Value = 41.77889500000001
$ws.Range("O5").Value = 0.1894184303222152
$ws.Range("P5").Value = 0.1894184303222152
$ws.Range("Q5").Value = 154.8990690182434
$ws.Range("R5").Value = 1394.09162116419
$ws.Range("S5").Value = 0.04640093523825884
$ws.Range("T5").Value = 0.04640093523825883
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.220714
$ws.Range("H6").Value = 54.662142
$ws.Range("I6").Value = 0.4012885600454987
$ws.Range("J6").Value = 0.4012885600454988
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.127188333333333
$ws.Range("N6").Value = 12.381565
$ws.Range("O6").Value = 0.0561359176022362
$ws.Range("P6").Value = 0.05613591760223619
$ws.Range("Q6").Value = 75.20031824580333
$ws.Range("R6").Value = 676.8028642122299
$ws.Range("S6").Value = 0.02252670154143413
$ws.Range("T6").Value = 0.02252670154143413
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.220714
$ws.Range("H7").Value = 54.662142
$ws.Range("I7").Value = 0.4012885600454987
$ws.Range("J7").Value = 0.4012885600454988
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 47.24901333333333
$ws.Range("N7").Value = 141.74704
$ws.Range("O7").Value = 0.6426570597336346
$ws.Range("P7").Value = 0.6426570597336345
$ws.Range("Q7").Value = 860.9107587288534
$ws.Range("R7").Value = 7748.19682855968
$ws.Range("S7").Value = 0.2578909261035843
$ws.Range("T7").Value = 0.2578909261035843
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.220714
$ws.Range("H8").Value = 54.662142
$ws.Range("I8").Value = 0.4012885600454987
$ws.Range("J8").Value = 0.4012885600454988
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.218847999999999
$ws.Range("N8").Value = 24.656544
$ws.Range("O8").Value = 0.1117885923419141
$ws.Range("P8").Value = 0.1117885923419141
$ws.Range("Q8").Value = 149.753278817472
$ws.Range("R8").Value = 1347.779509357248
$ws.Range("S8").Value = 0.04485948325039997
$ws.Range("T8").Value = 0.04485948325039998
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.220714
$ws.Range("H9").Value = 54.662142
$ws.Range("I9").Value = 0.4012885600454987
$ws.Range("J9").Value = 0.4012885600454988
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.92629833333334
$ws.Range("N9").Value = 41.77889500000001
$ws.Range("O9").Value = 0.1894184303222152
$ws.Range("P9").Value = 0.1894184303222152
$ws.Range("Q9").Value = 253.7470990103434
$ws.Range("R9").Value = 2283.723891093091
$ws.Range("S9").Value = 0.07601144915008039
$ws.Range("T9").Value = 0.07601144915008039
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1189986666666667
$ws.Range("H10").Value = 0.356996
$ws.Range("I10").Value = 0.002620797603979787
$ws.Range("J10").Value = 0.002620797603979787
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.127188333333333
$ws.Range("N10").Value = 12.381565
$ws.Range("O10").Value = 0.0561359176022362
$ws.Range("P10").Value = 0.05613591760223619
$ws.Range("Q10").Value = 0.4911299087488888
$ws.Range("R10").Value = 4.420169178739999
$ws.Range("S10").Value = 0.0001471208783491474
$ws.Range("T10").Value = 0.0001471208783491473
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1189986666666667
$ws.Range("H11").Value = 0.356996
$ws.Range("I11").Value = 0.002620797603979787
$ws.Range("J11").Value = 0.002620797603979787
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 47.24901333333333
$ws.Range("N11").Value = 141.74704
$ws.Range("O11").Value = 0.6426570597336346
$ws.Range("P11").Value = 0.6426570597336345
$ws.Range("Q11").Value = 5.622569587982222
$ws.Range("R11").Value = 50.60312629184
$ws.Range("S11").Value = 0.001684274082330604
$ws.Range("T11").Value = 0.001684274082330604
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1189986666666667
$ws.Range("H12").Value = 0.356996
$ws.Range("I12").Value = 0.002620797603979787
$ws.Range("J12").Value = 0.002620797603979787
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 8.218847999999999
$ws.Range("N12").Value = 24.656544
$ws.Range("O12").Value = 0.1117885923419141
$ws.Range("P12").Value = 0.1117885923419141
$ws.Range("Q12").Value = 0.9780319535359998
$ws.Range("R12").Value = 8.802287581824
$ws.Range("S12").Value = 0.0002929752749619616
$ws.Range("T12").Value = 0.0002929752749619616
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1189986666666667
$ws.Range("H13").Value = 0.356996
$ws.Range("I13").Value = 0.002620797603979787
$ws.Range("J13").Value = 0.002620797603979787
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.92629833333334
$ws.Range("N13").Value = 41.77889500000001
$ws.Range("O13").Value = 0.1894184303222152
$ws.Range("P13").Value = 0.1894184303222152
$ws.Range("Q13").Value = 1.657210933268889
$ws.Range("R13").Value = 14.91489839942
$ws.Range("S13").Value = 0.0004964273683380738
$ws.Range("T13").Value = 0.0004964273683380737
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.943029
$ws.Range("H14").Value = 47.829087
$ws.Range("I14").Value = 0.3511253812651704
$ws.Range("J14").Value = 0.3511253812651704
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.127188333333333
$ws.Range("N14").Value = 12.381565
$ws.Range("O14").Value = 0.0561359176022362
$ws.Range("P14").Value = 0.05613591760223619
$ws.Range("Q14").Value = 65.79988328679499
$ws.Range("R14").Value = 592.1989495811549
$ws.Range("S14").Value = 0.01971074547075537
$ws.Range("T14").Value = 0.01971074547075537
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.943029
$ws.Range("H15").Value = 47.829087
$ws.Range("I15").Value = 0.3511253812651704
$ws.Range("J15").Value = 0.3511253812651704
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 47.24901333333333
$ws.Range("N15").Value = 141.74704
$ws.Range("O15").Value = 0.6426570597336346
$ws.Range("P15").Value = 0.6426570597336345
$ws.Range("Q15").Value = 753.29238979472
$ws.Range("R15").Value = 6779.63150815248
$ws.Range("S15").Value = 0.2256532051217258
$ws.Range("T15").Value = 0.2256532051217258
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.943029
$ws.Range("H16").Value = 47.829087
$ws.Range("I16").Value = 0.3511253812651704
$ws.Range("J16").Value = 0.3511253812651704
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 8.218847999999999
$ws.Range("N16").Value = 24.656544
$ws.Range("O16").Value = 0.1117885923419141
$ws.Range("P16").Value = 0.1117885923419141
$ws.Range("Q16").Value = 131.033332010592
$ws.Range("R16").Value = 1179.299988095328
$ws.Range("S16").Value = 0.03925181210715129
$ws.Range("T16").Value = 0.03925181210715129
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.943029
$ws.Range("H17").Value = 47.829087
$ws.Range("I17").Value = 0.3511253812651704
$ws.Range("J17").Value = 0.3511253812651704
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.92629833333334
$ws.Range("N17").Value = 41.77889500000001
$ws.Range("O17").Value = 0.1894184303222152
$ws.Range("P17").Value = 0.1894184303222152
$ws.Range("Q17").Value = 222.027378190985
$ws.Range("R17").Value = 1998.246403718865
$ws.Range("S17").Value = 0.06650961856553793
$ws.Range("T17").Value = 0.06650961856553793
